# Adds the 18-sep column (CS) to "Prix Spot", and a new
# 2025-09-16 row (94) to both "Gaz" and "CO2" sheets,
# matching the daily automatic EPEX Spot data refresh.
$wb = $excel.ActiveWorkbook

# --- Sheet "Prix Spot": add column CS ("18-sep") ---
$ws1 = $wb.Worksheets.Item("Prix Spot")

# Copy the formatting of the previous header cell (CR1) so the new
# header (CS1) keeps the bold / centered / bordered style.
$ws1.Cells.Item(1,96).Copy($ws1.Cells.Item(1,97))
$ws1.Cells.Item(1,97).Value = "18-sep"

$ws1.Cells.Item(2,97).Value = 27.6
$ws1.Cells.Item(3,97).Value = 27.02
$ws1.Cells.Item(4,97).Value = 27.4
$ws1.Cells.Item(5,97).Value = 15.61
$ws1.Cells.Item(6,97).Value = 17.23
$ws1.Cells.Item(7,97).Value = 25.65
$ws1.Cells.Item(8,97).Value = 40.7
$ws1.Cells.Item(9,97).Value = 100.05
$ws1.Cells.Item(10,97).Value = 102.63
$ws1.Cells.Item(11,97).Value = 79.09
$ws1.Cells.Item(12,97).Value = 56.4
$ws1.Cells.Item(13,97).Value = 5.87
$ws1.Cells.Item(14,97).Value = 0.2
$ws1.Cells.Item(15,97).Value = 0
$ws1.Cells.Item(16,97).Value = 0
$ws1.Cells.Item(17,97).Value = 0.52
$ws1.Cells.Item(18,97).Value = 15.31
$ws1.Cells.Item(19,97).Value = 45.8
$ws1.Cells.Item(20,97).Value = 96.86
$ws1.Cells.Item(21,97).Value = 124.7
$ws1.Cells.Item(22,97).Value = 127.5
$ws1.Cells.Item(23,97).Value = 110
$ws1.Cells.Item(24,97).Value = 99.09
$ws1.Cells.Item(25,97).Value = 88.2

# --- Sheet "Gaz": append row 94 (2025-09-16) ---
$ws2 = $wb.Worksheets.Item("Gaz")

# Force the date column to stay plain text (matching the existing
# "YYYY-MM-DD" string cells) instead of Excel auto-converting it to a
# date serial number, then drop the temporary Text format again so the
# cell keeps the default (unstyled) look of the rest of the column.
$ws2.Cells.Item(94,1).NumberFormat = "@"
$ws2.Cells.Item(94,1).Value = "2025-09-16"
$ws2.Cells.Item(94,1).ClearFormats()
$ws2.Cells.Item(94,2).Value = 31.925

# --- Sheet "CO2": append row 94 (2025-09-16) ---
$ws3 = $wb.Worksheets.Item("CO2")
$ws3.Cells.Item(94,1).NumberFormat = "@"
$ws3.Cells.Item(94,1).Value = "2025-09-16"
$ws3.Cells.Item(94,1).ClearFormats()
$ws3.Cells.Item(94,2).Value = 77.29000000000001
